# Daily "remaining days" refresh for the shop-tracking sheet.
#
# Columns: D = total days (总天), E = remaining days (剩余), F = start date
# as an 8-digit yyyyMMdd integer (开始时间).
#
# For every data row the prior snapshot satisfies:
#     E = D - (today_before - F)
# i.e. "today_before" is implied consistently by every row's D/E/F triple.
# We recover that implied date straight from the sheet (majority vote,
# so a row with a corrupted F is simply ignored), advance it by one day
# (this script represents the next day's run), and recompute E for every
# row against "today_after". Whenever a row would run out (new E < 1) the
# stock is treated as replenished today: E resets to D and F becomes
# today_after (yyyyMMdd).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

function ConvertTo-DateFromYyyymmdd($n) {
    if ($n -eq $null) { return $null }
    $s = [string]$n
    if ($s.Length -ne 8) { return $null }
    $y = [int]$s.Substring(0,4)
    $m = [int]$s.Substring(4,2)
    $d = [int]$s.Substring(6,2)
    if ($m -lt 1 -or $m -gt 12) { return $null }
    if ($d -lt 1 -or $d -gt 31) { return $null }
    try {
        return (Get-Date -Year $y -Month $m -Day $d -Hour 0 -Minute 0 -Second 0)
    } catch {
        return $null
    }
}

function Get-DayDiff($dateFrom, $dateTo) {
    # ($dateTo - $dateFrom) in whole days, via OADate (TimeSpan subtraction
    # is unreliable for DateTime operands here).
    $diff = $dateTo.ToOADate() - $dateFrom.ToOADate()
    return [int][math]::Round($diff)
}

# ---- Pass 1: recover "today" (before this run) from the existing rows ----
$votes = @{}
for ($r = 2; $r -le $lastRow; $r++) {
    $dVal = $ws.Cells.Item($r, 4).Value2
    $eVal = $ws.Cells.Item($r, 5).Value2
    $fVal = $ws.Cells.Item($r, 6).Value2
    if ($dVal -eq $null -or $eVal -eq $null -or $fVal -eq $null) { continue }

    $fDate = ConvertTo-DateFromYyyymmdd $fVal
    if ($fDate -eq $null) { continue }

    $elapsed = [int]$dVal - [int]$eVal
    $candidate = $fDate.AddDays($elapsed)
    $key = $candidate.ToString("yyyyMMdd")
    if ($votes.ContainsKey($key)) {
        $votes[$key] = $votes[$key] + 1
    } else {
        $votes[$key] = 1
    }
}

$bestKey = $null
$bestCount = -1
foreach ($k in $votes.Keys) {
    if ($votes[$k] -gt $bestCount) {
        $bestCount = $votes[$k]
        $bestKey = $k
    }
}

$todayBefore = ConvertTo-DateFromYyyymmdd $bestKey
$todayAfter = $todayBefore.AddDays(1)
$todayAfterStr = $todayAfter.ToString("yyyyMMdd")
$todayAfterInt = [int]$todayAfterStr

# ---- Pass 2: recompute E (and F when stock is replenished) per row ----
for ($r = 2; $r -le $lastRow; $r++) {
    $dVal = $ws.Cells.Item($r, 4).Value2
    $eVal = $ws.Cells.Item($r, 5).Value2
    $fVal = $ws.Cells.Item($r, 6).Value2
    if ($dVal -eq $null -or $eVal -eq $null -or $fVal -eq $null) { continue }

    $fDate = ConvertTo-DateFromYyyymmdd $fVal
    if ($fDate -eq $null) { continue }

    $dInt = [int]$dVal

    $elapsedAfter = Get-DayDiff $fDate $todayAfter
    $newE = $dInt - $elapsedAfter

    if ($newE -lt 1) {
        $ws.Cells.Item($r, 5).Value = $dInt
        $ws.Cells.Item($r, 6).Value = $todayAfterInt
    } else {
        $ws.Cells.Item($r, 5).Value = $newE
    }
}
